$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = [double]"1.002365"
$ws.Range("H2").Value = [double]"3.007095"
$ws.Range("I2").Value = [double]"0.00427078073065012"
$ws.Range("J2").Value = [double]"0.00427078073065012"
$ws.Range("M2").Value = [double]"0.3331066666666667"
$ws.Range("N2").Value = [double]"0.99932"
$ws.Range("O2").Value = [double]"0.002125805913843485"
$ws.Range("P2").Value = [double]"0.002125805913843485"
$ws.Range("Q2").Value = [double]"0.3338944639333333"
$ws.Range("R2").Value = [double]"3.0050501754"
$ws.Range("S2").Value = [double]"9.078850933944825E-06"
$ws.Range("T2").Value = [double]"9.078850933944825E-06"
$ws.Range("G3").Value = [double]"1.002365"
$ws.Range("H3").Value = [double]"3.007095"
$ws.Range("I3").Value = [double]"0.00427078073065012"
$ws.Range("J3").Value = [double]"0.00427078073065012"
$ws.Range("O3").Value = [double]"0.0008775937418887864"
$ws.Range("P3").Value = [double]"0.0008775937418887864"
$ws.Range("Q3").Value = [double]"0.13784122534"
$ws.Range("R3").Value = [double]"1.24057102806"
$ws.Range("S3").Value = [double]"3.748010442197764E-06"
$ws.Range("T3").Value = [double]"3.748010442197764E-06"
$ws.Range("G4").Value = [double]"1.002365"
$ws.Range("H4").Value = [double]"3.007095"
$ws.Range("I4").Value = [double]"0.00427078073065012"
$ws.Range("J4").Value = [double]"0.00427078073065012"
$ws.Range("M4").Value = [double]"91.40156066666667"
$ws.Range("N4").Value = [double]"274.204682"
$ws.Range("O4").Value = [double]"0.5833025803538128"
$ws.Range("P4").Value = [double]"0.5833025803538128"
$ws.Range("Q4").Value = [double]"91.61772535764334"
$ws.Range("R4").Value = [double]"824.5595282187899"
$ws.Range("S4").Value = [double]"0.002491157420313557"
$ws.Range("T4").Value = [double]"0.002491157420313557"
$ws.Range("G5").Value = [double]"1.002365"
$ws.Range("H5").Value = [double]"3.007095"
$ws.Range("I5").Value = [double]"0.00427078073065012"
$ws.Range("J5").Value = [double]"0.00427078073065012"
$ws.Range("M5").Value = [double]"0.5759770000000001"
$ws.Range("N5").Value = [double]"1.727931"
$ws.Range("O5").Value = [double]"0.00367574544541637"
$ws.Range("P5").Value = [double]"0.00367574544541637"
$ws.Range("Q5").Value = [double]"0.5773391856050001"
$ws.Range("R5").Value = [double]"5.196052670445"
$ws.Range("S5").Value = [double]"1.569830281905918E-05"
$ws.Range("T5").Value = [double]"1.569830281905918E-05"
$ws.Range("G6").Value = [double]"1.002365"
$ws.Range("H6").Value = [double]"3.007095"
$ws.Range("I6").Value = [double]"0.00427078073065012"
$ws.Range("J6").Value = [double]"0.00427078073065012"
$ws.Range("M6").Value = [double]"64.24849033333334"
$ws.Range("N6").Value = [double]"192.745471"
$ws.Range("O6").Value = [double]"0.4100182745450386"
$ws.Range("P6").Value = [double]"0.4100182745450385"
$ws.Range("Q6").Value = [double]"64.40043801297166"
$ws.Range("R6").Value = [double]"579.603942116745"
$ws.Range("S6").Value = [double]"0.001751098146141362"
$ws.Range("T6").Value = [double]"0.001751098146141361"
$ws.Range("I7").Value = [double]"0.00206557659722455"
$ws.Range("J7").Value = [double]"0.002065576597224551"
$ws.Range("M7").Value = [double]"0.3331066666666667"
$ws.Range("N7").Value = [double]"0.99932"
$ws.Range("O7").Value = [double]"0.002125805913843485"
$ws.Range("P7").Value = [double]"0.002125805913843485"
$ws.Range("Q7").Value = [double]"0.16148911268"
$ws.Range("R7").Value = [double]"1.45340201412"
$ws.Range("S7").Value = [double]"4.391014945876651E-06"
$ws.Range("T7").Value = [double]"4.391014945876651E-06"
$ws.Range("I8").Value = [double]"0.00206557659722455"
$ws.Range("J8").Value = [double]"0.002065576597224551"
$ws.Range("O8").Value = [double]"0.0008775937418887864"
$ws.Range("P8").Value = [double]"0.0008775937418887864"
$ws.Range("S8").Value = [double]"1.8127370951162E-06"
$ws.Range("T8").Value = [double]"1.8127370951162E-06"
$ws.Range("I9").Value = [double]"0.00206557659722455"
$ws.Range("J9").Value = [double]"0.002065576597224551"
$ws.Range("M9").Value = [double]"91.40156066666667"
$ws.Range("N9").Value = [double]"274.204682"
$ws.Range("O9").Value = [double]"0.5833025803538128"
$ws.Range("P9").Value = [double]"0.5833025803538128"
$ws.Range("Q9").Value = [double]"44.311202406518"
$ws.Range("R9").Value = [double]"398.800821658662"
$ws.Range("S9").Value = [double]"0.001204856159079528"
$ws.Range("T9").Value = [double]"0.001204856159079529"
$ws.Range("I10").Value = [double]"0.00206557659722455"
$ws.Range("J10").Value = [double]"0.002065576597224551"
$ws.Range("M10").Value = [double]"0.5759770000000001"
$ws.Range("N10").Value = [double]"1.727931"
$ws.Range("O10").Value = [double]"0.00367574544541637"
$ws.Range("P10").Value = [double]"0.00367574544541637"
$ws.Range("Q10").Value = [double]"0.279231921669"
$ws.Range("R10").Value = [double]"2.513087295021"
$ws.Range("S10").Value = [double]"7.592533769406785E-06"
$ws.Range("T10").Value = [double]"7.592533769406786E-06"
$ws.Range("I11").Value = [double]"0.00206557659722455"
$ws.Range("J11").Value = [double]"0.002065576597224551"
$ws.Range("M11").Value = [double]"64.24849033333334"
$ws.Range("N11").Value = [double]"192.745471"
$ws.Range("O11").Value = [double]"0.4100182745450386"
$ws.Range("P11").Value = [double]"0.4100182745450385"
$ws.Range("Q11").Value = [double]"31.147475368129"
$ws.Range("R11").Value = [double]"280.327278313161"
$ws.Range("S11").Value = [double]"0.0008469241523346222"
$ws.Range("T11").Value = [double]"0.0008469241523346222"
$ws.Range("G12").Value = [double]"135.9134216666667"
$ws.Range("H12").Value = [double]"407.740265"
$ws.Range("I12").Value = [double]"0.5790868818152315"
$ws.Range("J12").Value = [double]"0.5790868818152316"
$ws.Range("M12").Value = [double]"0.3331066666666667"
$ws.Range("N12").Value = [double]"0.99932"
$ws.Range("O12").Value = [double]"0.002125805913843485"
$ws.Range("P12").Value = [double]"0.002125805913843485"
$ws.Range("Q12").Value = [double]"45.27366684664445"
$ws.Range("R12").Value = [double]"407.4630016198"
$ws.Range("S12").Value = [double]"0.001231026317992002"
$ws.Range("T12").Value = [double]"0.001231026317992003"
$ws.Range("G13").Value = [double]"135.9134216666667"
$ws.Range("H13").Value = [double]"407.740265"
$ws.Range("I13").Value = [double]"0.5790868818152315"
$ws.Range("J13").Value = [double]"0.5790868818152316"
$ws.Range("O13").Value = [double]"0.0008775937418887864"
$ws.Range("P13").Value = [double]"0.0008775937418887864"
$ws.Range("Q13").Value = [double]"18.69027009391333"
$ws.Range("R13").Value = [double]"168.21243084522"
$ws.Range("S13").Value = [double]"0.0005082030234909384"
$ws.Range("T13").Value = [double]"0.0005082030234909385"
$ws.Range("G14").Value = [double]"135.9134216666667"
$ws.Range("H14").Value = [double]"407.740265"
$ws.Range("I14").Value = [double]"0.5790868818152315"
$ws.Range("J14").Value = [double]"0.5790868818152316"
$ws.Range("M14").Value = [double]"91.40156066666667"
$ws.Range("N14").Value = [double]"274.204682"
$ws.Range("O14").Value = [double]"0.5833025803538128"
$ws.Range("P14").Value = [double]"0.5833025803538128"
$ws.Range("Q14").Value = [double]"12422.69885588008"
$ws.Range("R14").Value = [double]"111804.2897029207"
$ws.Range("S14").Value = [double]"0.3377828724118679"
$ws.Range("T14").Value = [double]"0.337782872411868"
$ws.Range("G15").Value = [double]"135.9134216666667"
$ws.Range("H15").Value = [double]"407.740265"
$ws.Range("I15").Value = [double]"0.5790868818152315"
$ws.Range("J15").Value = [double]"0.5790868818152316"
$ws.Range("M15").Value = [double]"0.5759770000000001"
$ws.Range("N15").Value = [double]"1.727931"
$ws.Range("O15").Value = [double]"0.00367574544541637"
$ws.Range("P15").Value = [double]"0.00367574544541637"
$ws.Range("Q15").Value = [double]"78.28300487130167"
$ws.Range("R15").Value = [double]"704.5470438417151"
$ws.Range("S15").Value = [double]"0.002128575968332705"
$ws.Range("T15").Value = [double]"0.002128575968332705"
$ws.Range("G16").Value = [double]"135.9134216666667"
$ws.Range("H16").Value = [double]"407.740265"
$ws.Range("I16").Value = [double]"0.5790868818152315"
$ws.Range("J16").Value = [double]"0.5790868818152316"
$ws.Range("M16").Value = [double]"64.24849033333334"
$ws.Range("N16").Value = [double]"192.745471"
$ws.Range("O16").Value = [double]"0.4100182745450386"
$ws.Range("P16").Value = [double]"0.4100182745450385"
$ws.Range("Q16").Value = [double]"8732.232158121091"
$ws.Range("R16").Value = [double]"78590.08942308983"
$ws.Range("S16").Value = [double]"0.2374362040935479"
$ws.Range("T16").Value = [double]"0.2374362040935479"
$ws.Range("G17").Value = [double]"0.06627866666666667"
$ws.Range("H17").Value = [double]"0.198836"
$ws.Range("I17").Value = [double]"0.0002823937911371431"
$ws.Range("J17").Value = [double]"0.0002823937911371431"
$ws.Range("M17").Value = [double]"0.3331066666666667"
$ws.Range("N17").Value = [double]"0.99932"
$ws.Range("O17").Value = [double]"0.002125805913843485"
$ws.Range("P17").Value = [double]"0.002125805913843485"
$ws.Range("Q17").Value = [double]"0.02207786572444444"
$ws.Range("R17").Value = [double]"0.19870079152"
$ws.Range("S17").Value = [double]"6.003143912320206E-07"
$ws.Range("T17").Value = [double]"6.003143912320207E-07"
$ws.Range("G18").Value = [double]"0.06627866666666667"
$ws.Range("H18").Value = [double]"0.198836"
$ws.Range("I18").Value = [double]"0.0002823937911371431"
$ws.Range("J18").Value = [double]"0.0002823937911371431"
$ws.Range("O18").Value = [double]"0.0008775937418887864"
$ws.Range("P18").Value = [double]"0.0008775937418887864"
$ws.Range("Q18").Value = [double]"0.009114377125333334"
$ws.Range("R18").Value = [double]"0.082029394128"
$ws.Range("S18").Value = [double]"2.478270238502058E-07"
$ws.Range("T18").Value = [double]"2.478270238502058E-07"
$ws.Range("G19").Value = [double]"0.06627866666666667"
$ws.Range("H19").Value = [double]"0.198836"
$ws.Range("I19").Value = [double]"0.0002823937911371431"
$ws.Range("J19").Value = [double]"0.0002823937911371431"
$ws.Range("M19").Value = [double]"91.40156066666667"
$ws.Range("N19").Value = [double]"274.204682"
$ws.Range("O19").Value = [double]"0.5833025803538128"
$ws.Range("P19").Value = [double]"0.5833025803538128"
$ws.Range("Q19").Value = [double]"6.057973572239111"
$ws.Range("R19").Value = [double]"54.521762150152"
$ws.Range("S19").Value = [double]"0.0001647210270461912"
$ws.Range("T19").Value = [double]"0.0001647210270461913"
$ws.Range("G20").Value = [double]"0.06627866666666667"
$ws.Range("H20").Value = [double]"0.198836"
$ws.Range("I20").Value = [double]"0.0002823937911371431"
$ws.Range("J20").Value = [double]"0.0002823937911371431"
$ws.Range("M20").Value = [double]"0.5759770000000001"
$ws.Range("N20").Value = [double]"1.727931"
$ws.Range("O20").Value = [double]"0.00367574544541637"
$ws.Range("P20").Value = [double]"0.00367574544541637"
$ws.Range("Q20").Value = [double]"0.03817498759066667"
$ws.Range("R20").Value = [double]"0.343574888316"
$ws.Range("S20").Value = [double]"1.038007691586216E-06"
$ws.Range("T20").Value = [double]"1.038007691586216E-06"
$ws.Range("G21").Value = [double]"0.06627866666666667"
$ws.Range("H21").Value = [double]"0.198836"
$ws.Range("I21").Value = [double]"0.0002823937911371431"
$ws.Range("J21").Value = [double]"0.0002823937911371431"
$ws.Range("M21").Value = [double]"64.24849033333334"
$ws.Range("N21").Value = [double]"192.745471"
$ws.Range("O21").Value = [double]"0.4100182745450386"
$ws.Range("P21").Value = [double]"0.4100182745450385"
$ws.Range("Q21").Value = [double]"4.258304274639555"
$ws.Range("R21").Value = [double]"38.32473847175601"
$ws.Range("S21").Value = [double]"0.0001157866149842834"
$ws.Range("T21").Value = [double]"0.0001157866149842834"
$ws.Range("G22").Value = [double]"97.236126"
$ws.Range("H22").Value = [double]"291.708378"
$ws.Range("I22").Value = [double]"0.4142943670657566"
$ws.Range("J22").Value = [double]"0.4142943670657567"
$ws.Range("M22").Value = [double]"0.3331066666666667"
$ws.Range("N22").Value = [double]"0.99932"
$ws.Range("O22").Value = [double]"0.002125805913843485"
$ws.Range("P22").Value = [double]"0.002125805913843485"
$ws.Range("Q22").Value = [double]"32.39000181144"
$ws.Range("R22").Value = [double]"291.51001630296"
$ws.Range("S22").Value = [double]"0.0008807094155804289"
$ws.Range("T22").Value = [double]"0.000880709415580429"
$ws.Range("G23").Value = [double]"97.236126"
$ws.Range("H23").Value = [double]"291.708378"
$ws.Range("I23").Value = [double]"0.4142943670657566"
$ws.Range("J23").Value = [double]"0.4142943670657567"
$ws.Range("O23").Value = [double]"0.0008775937418887864"
$ws.Range("P23").Value = [double]"0.0008775937418887864"
$ws.Range("Q23").Value = [double]"13.371523103016"
$ws.Range("R23").Value = [double]"120.343707927144"
$ws.Range("S23").Value = [double]"0.0003635821438366837"
$ws.Range("T23").Value = [double]"0.0003635821438366838"
$ws.Range("G24").Value = [double]"97.236126"
$ws.Range("H24").Value = [double]"291.708378"
$ws.Range("I24").Value = [double]"0.4142943670657566"
$ws.Range("J24").Value = [double]"0.4142943670657567"
$ws.Range("M24").Value = [double]"91.40156066666667"
$ws.Range("N24").Value = [double]"274.204682"
$ws.Range("O24").Value = [double]"0.5833025803538128"
$ws.Range("P24").Value = [double]"0.5833025803538128"
$ws.Range("Q24").Value = [double]"8887.533669580644"
$ws.Range("R24").Value = [double]"79987.80302622578"
$ws.Range("S24").Value = [double]"0.2416589733355055"
$ws.Range("T24").Value = [double]"0.2416589733355055"
$ws.Range("G25").Value = [double]"97.236126"
$ws.Range("H25").Value = [double]"291.708378"
$ws.Range("I25").Value = [double]"0.4142943670657566"
$ws.Range("J25").Value = [double]"0.4142943670657567"
$ws.Range("M25").Value = [double]"0.5759770000000001"
$ws.Range("N25").Value = [double]"1.727931"
$ws.Range("O25").Value = [double]"0.00367574544541637"
$ws.Range("P25").Value = [double]"0.00367574544541637"
$ws.Range("Q25").Value = [double]"56.00577214510201"
$ws.Range("R25").Value = [double]"504.051949305918"
$ws.Range("S25").Value = [double]"0.001522840632803613"
$ws.Range("T25").Value = [double]"0.001522840632803613"
$ws.Range("G26").Value = [double]"97.236126"
$ws.Range("H26").Value = [double]"291.708378"
$ws.Range("I26").Value = [double]"0.4142943670657566"
$ws.Range("J26").Value = [double]"0.4142943670657567"
$ws.Range("M26").Value = [double]"64.24849033333334"
$ws.Range("N26").Value = [double]"192.745471"
$ws.Range("O26").Value = [double]"0.4100182745450386"
$ws.Range("P26").Value = [double]"0.4100182745450385"
$ws.Range("Q26").Value = [double]"6247.274301361782"
$ws.Range("R26").Value = [double]"56225.46871225604"
$ws.Range("S26").Value = [double]"0.1698682615380304"
$ws.Range("T26").Value = [double]"0.1698682615380304"
